$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 06:58:17"
$wsZhCn.Range("H2").Value = "2016-03-22 06:59:02"

# "de-de" sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 06:58:27"
$wsDeDe.Range("H2").Value = "2016-03-22 06:59:15"
